# Refresh cryptos price/volume snapshot (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.143.72'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '2.565.95'
$ws.Range('E3').Value = '  +0.96%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E5').Value = '  +3.07%  '
$ws.Range('D6').Value = "'147.66"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.65%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = "'0.600"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.43%  '
$ws.Range('E9').Value = '  +3.49%  '
$ws.Range('E10').Value = '  +0.77%  '
$ws.Range('E11').Value = '  +0.01%  '
$ws.Range('E12').Value = '  +1.51%  '
$ws.Range('D13').Value = "'27.40"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('D14').Value = '3.026.44'
$ws.Range('E14').Value = '  +0.98%  '
$ws.Range('D15').Value = '63.106.96'
$ws.Range('E15').Value = '  +0.42%  '
$ws.Range('E16').Value = '  +3.40%  '
$ws.Range('D17').Value = '2.566.81'
$ws.Range('E17').Value = '  +0.92%  '
$ws.Range('E18').Value = '  -0.58%  '
$ws.Range('D19').Value = "'344.16"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.26%  '
$ws.Range('E20').Value = '  +2.96%  '
$ws.Range('E21').Value = '  +2.30%  '
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('E23').Value = '  -3.84%  '
$ws.Range('D24').Value = "'66.82"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.43%  '
$ws.Range('D25').Value = '2.697.06'
$ws.Range('E25').Value = '  +0.88%  '
$ws.Range('E26').Value = '  +1.05%  '
$ws.Range('D27').Value = "'1.63"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.67%  '
$ws.Range('D28').Value = "'8.10"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +11.07%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').Value = "'1.00"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = "'8.46"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.89%  '
$ws.Range('E31').Value = '  -1.44%  '
$ws.Range('E32').Value = '  +9.35%  '
$ws.Range('D33').Value = '0.0₃0826'
$ws.Range('E33').Value = '  +1.65%  '
$ws.Range('D34').Value = "'464.72"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +14.30%  '
$ws.Range('E35').Value = '  +3.55%  '
$ws.Range('D36').Value = "'175.93"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.20%  '
$ws.Range('E37').Value = '  +2.33%  '
$ws.Range('D38').Value = "'19.22"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.43%  '
$ws.Range('D39').Value = "'4.56"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +4.79%  '
$ws.Range('E41').Value = '  -0.69%  '
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('D43').Value = "'151.34"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.03%  '
$ws.Range('E44').Value = '  +2.30%  '
$ws.Range('D45').Value = "'20.97"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.04%  '
$ws.Range('E46').Value = '  +5.90%  '
$ws.Range('E47').Value = '  +1.12%  '
$ws.Range('E48').Value = '  +1.75%  '
$ws.Range('E49').Value = '  +1.76%  '
$ws.Range('D50').Value = "'1.74"
$ws.Range('D50').ClearFormats()
$ws.Range('E51').Value = '  -0.09%  '
